# Append the latest EUR -> ARS quotation as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

# Column A/B hold date-looking / time-looking text ("2025-10-18", "15:18:24").
# Force the cell format to Text before assigning so Excel stores the literal
# string instead of auto-converting it to a date/time serial number, then
# restore the default "Normal" style so no extra number-format style is left
# attached to the cell (matches the rest of the sheet, which uses the
# default style throughout).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-18"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "15:18:24"
$ws.Cells.Item($row, 2).Style = "Normal"

# Column C is unambiguous text ("1.00 EUR = 1,703.9130") so no special
# handling is required.
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,703.9130"
